$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2407.8845
$ws.Range("J19").Value = 560.9231
$ws.Range("L19").Value = 560.9231
$ws.Range("N19").Value = -910.9231
$ws.Range("H116").Value = 103419.336
$ws.Range("I116").Value = 134212.88
$ws.Range("K116").Value = 134212.88
$ws.Range("M116").Value = -130770.88

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 434.1
$ws.Range("I2").Value = 371.22223
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 371.22223
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -258.22223
$ws.Range("N2").Value = -1226
$ws.Range("H61").Value = 3122.0408
$ws.Range("I61").Value = 2031.3572
$ws.Range("K61").Value = 2031.3572
$ws.Range("M61").Value = -1819.3572
$ws.Range("H63").Value = 4546.8276
$ws.Range("I63").Value = 4819.0835
$ws.Range("K63").Value = 4819.0835
$ws.Range("M63").Value = -4133.0835
$ws.Range("H66").Value = 4546.8276
$ws.Range("I66").Value = 4819.0835
$ws.Range("K66").Value = 24095.4175
$ws.Range("M66").Value = -20663.4175
$ws.Range("H116").Value = 434.1
$ws.Range("I116").Value = 371.22223
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 371.22223
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = 1922.77777
$ws.Range("N116").Value = -5588
$ws.Range("H132").Value = 4474.794
$ws.Range("I132").Value = 5381.8237
$ws.Range("K132").Value = 16145.4711
$ws.Range("M132").Value = -13615.4711
$ws.Range("H136").Value = 3122.0408
$ws.Range("I136").Value = 2031.3572
$ws.Range("K136").Value = 6094.071599999999
$ws.Range("M136").Value = -3544.071599999999

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 434.1
$ws.Range("I3").Value = 371.22223
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 371.22223
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -257.22223
$ws.Range("N3").Value = -1228
$ws.Range("H94").Value = 782.71875
$ws.Range("I94").Value = 652.3137
$ws.Range("J94").Value = 1294.3077
$ws.Range("K94").Value = 652.3137
$ws.Range("L94").Value = 1294.3077
$ws.Range("M94").Value = -201.3137
$ws.Range("N94").Value = -2196.3077

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2416.8333
$ws.Range("I16").Value = 3175.25
$ws.Range("K16").Value = 3175.25
$ws.Range("M16").Value = -2888.25
$ws.Range("H31").Value = 3684
$ws.Range("I31").Value = 2653.6296
$ws.Range("J31").Value = 5671.143
$ws.Range("K31").Value = 2653.6296
$ws.Range("L31").Value = 5671.143
$ws.Range("M31").Value = -2358.6296
$ws.Range("N31").Value = -6261.143
$ws.Range("H34").Value = 3684
$ws.Range("I34").Value = 2653.6296
$ws.Range("J34").Value = 5671.143
$ws.Range("K34").Value = 2653.6296
$ws.Range("L34").Value = 5671.143
$ws.Range("M34").Value = -2451.6296
$ws.Range("N34").Value = -6075.143
$ws.Range("H58").Value = 2070.2666
$ws.Range("I58").Value = 1505.1428
$ws.Range("J58").Value = 3388.889
$ws.Range("K58").Value = 1505.1428
$ws.Range("L58").Value = 3388.889
$ws.Range("M58").Value = -1302.1428
$ws.Range("N58").Value = -3794.889
$ws.Range("H99").Value = 50370.145
$ws.Range("J99").Value = 2987.7144
$ws.Range("L99").Value = 2987.7144
$ws.Range("N99").Value = -5983.7144
$ws.Range("H113").Value = 2416.8333
$ws.Range("I113").Value = 3175.25
$ws.Range("K113").Value = 3175.25
$ws.Range("M113").Value = -1005.25
$ws.Range("H122").Value = 1413.8864
$ws.Range("I122").Value = 1534.8889
$ws.Range("K122").Value = 4604.6667
$ws.Range("M122").Value = -2154.6667
$ws.Range("H126").Value = 50370.145
$ws.Range("J126").Value = 2987.7144
$ws.Range("L126").Value = 8963.143199999999
$ws.Range("N126").Value = -13903.1432
$ws.Range("H136").Value = 2070.2666
$ws.Range("I136").Value = 1505.1428
$ws.Range("J136").Value = 3388.889
$ws.Range("K136").Value = 4515.428400000001
$ws.Range("L136").Value = 10166.667
$ws.Range("M136").Value = -1965.428400000001
$ws.Range("N136").Value = -15266.667

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 966.6667
$ws.Range("I92").Value = 900
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 2700
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = -1452
$ws.Range("N92").Value = -5496
$ws.Range("H97").Value = 420.13333
$ws.Range("I97").Value = 123.818184
$ws.Range("J97").Value = 1235
$ws.Range("K97").Value = 371.454552
$ws.Range("L97").Value = 3705
$ws.Range("M97").Value = 124.545448
$ws.Range("N97").Value = -4697
$ws.Range("H131").Value = 2160.1702
$ws.Range("I131").Value = 4212.25
$ws.Range("J131").Value = 1739.2307
$ws.Range("K131").Value = 12636.75
$ws.Range("L131").Value = 5217.6921
$ws.Range("M131").Value = -7596.75
$ws.Range("N131").Value = -15297.6921
$ws.Range("H132").Value = 4388.839
$ws.Range("I132").Value = 3278.2222
$ws.Range("K132").Value = 29503.9998
$ws.Range("M132").Value = -26973.9998

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1028.6923
$ws.Range("I122").Value = 988.05
$ws.Range("J122").Value = 1164.1666
$ws.Range("K122").Value = 2964.15
$ws.Range("L122").Value = 3492.4998
$ws.Range("M122").Value = -514.1499999999996
$ws.Range("N122").Value = -8392.4998

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2056.3225
$ws.Range("I7").Value = 1891.5264
$ws.Range("J7").Value = 2317.25
$ws.Range("K7").Value = 1891.5264
$ws.Range("L7").Value = 2317.25
$ws.Range("M7").Value = -1779.5264
$ws.Range("N7").Value = -2541.25
$ws.Range("H126").Value = 2056.3225
$ws.Range("I126").Value = 1891.5264
$ws.Range("J126").Value = 2317.25
$ws.Range("K126").Value = 5674.5792
$ws.Range("L126").Value = 6951.75
$ws.Range("M126").Value = -3204.5792
$ws.Range("N126").Value = -11891.75
$ws.Range("H132").Value = 7777.054
$ws.Range("I132").Value = 2627.7222
$ws.Range("J132").Value = 12655.368
$ws.Range("K132").Value = 7883.1666
$ws.Range("L132").Value = 37966.104
$ws.Range("M132").Value = -5353.1666
$ws.Range("N132").Value = -43026.104

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15151.507
$ws.Range("I132").Value = 21167.408
$ws.Range("J132").Value = 2869.0417
$ws.Range("K132").Value = 63502.224
$ws.Range("L132").Value = 8607.125100000001
$ws.Range("M132").Value = -60972.224
$ws.Range("N132").Value = -13667.1251
$ws.Range("H136").Value = 16668415
$ws.Range("I136").Value = 31251406
$ws.Range("J136").Value = 2139.3928
$ws.Range("K136").Value = 93754218
$ws.Range("L136").Value = 6418.178400000001
$ws.Range("M136").Value = -93751668
$ws.Range("N136").Value = -11518.1784
